# Added 4wk low sales check: refresh the forecast summary figures to reflect
# the high-volume season trend and the new low-sales threshold logic.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: rows 2-17 ---
# Columns: D=MyForecast, G=Trend, H=Inventory Coverage,
#          I=Stockout Risk, J=Reorder Urgency, L=Seasonality Index

$rows = @(
    @{ Row = 2;  D = 111; G = "High Volume Season"; H = 14.13; I = "Low";  J = "Normal"; L = 0.97 },
    @{ Row = 3;  D = 149; G = "High Volume Season"; H = 9.81;  I = "Low";  J = "Normal"; L = 1 },
    @{ Row = 4;  D = 195; G = "High Volume Season"; H = 6.75;  I = "Low";  J = "Normal"; L = 0.84 },
    @{ Row = 5;  D = 218; G = "High Volume Season"; H = 5.14;  I = "Low";  J = "Normal"; L = 1.17 },
    @{ Row = 6;  D = 191; G = "High Volume Season"; H = 4.72;  I = "Low";  J = "Normal"; L = 1.09 },
    @{ Row = 7;  D = 122; G = "High Volume Season"; H = 5.83;  I = "Low";  J = "Normal"; L = 0.84 },
    @{ Row = 8;  D = 69;  G = "High Volume Season"; H = 8.51;  I = "Low";  J = "Normal"; L = 1.2 },
    @{ Row = 9;  D = 69;  G = "High Volume Season"; H = 7.51;  I = "Low";  J = "Normal"; L = 0.92 },
    @{ Row = 10; D = 111; G = "High Volume Season"; H = 4.04;  I = "Low";  J = "Normal"; L = 1.14 },
    @{ Row = 11; D = 167; G = "High Volume Season"; H = 2.04;  I = "Low";  J = "Normal"; L = 0.98 },
    @{ Row = 12; D = 167; G = "High Volume Season"; H = 1.03;  I = "Low";  J = "Normal"; L = 1.04 },
    @{ Row = 13; D = 118; G = "High Volume Season"; H = 0.04;  I = "High"; J = "Urgent"; L = 0.86 },
    @{ Row = 14; D = 78;  G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 0.8 },
    @{ Row = 15; D = 96;  G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 1.12 },
    @{ Row = 16; D = 153; G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 1.1 },
    @{ Row = 17; D = 183; G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 0.97 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $wsForecast.Cells.Item($row, 4).Value  = $r.D   # D: MyForecast
    $wsForecast.Cells.Item($row, 7).Value  = $r.G   # G: Trend
    $wsForecast.Cells.Item($row, 8).Value  = $r.H   # H: Inventory Coverage
    $wsForecast.Cells.Item($row, 9).Value  = $r.I   # I: Stockout Risk
    $wsForecast.Cells.Item($row, 10).Value = $r.J   # J: Reorder Urgency
    $wsForecast.Cells.Item($row, 12).Value = $r.L   # L: Seasonality Index
}

# --- Summary sheet: updated totals from the revised forecast ---
# (Source cells are stored as text, e.g. "1895"; the leading apostrophe
# keeps the numeric-looking value typed as text instead of a number,
# matching the original inline-string cell type.)
$wsSummary.Cells.Item(9, 2).Value  = "'2204"  # Total Forecast (16 Weeks)
$wsSummary.Cells.Item(10, 2).Value = "'1128"  # Total Forecast (8 Weeks)
$wsSummary.Cells.Item(11, 2).Value = "'675"   # Total Forecast (4 Weeks)
$wsSummary.Cells.Item(12, 2).Value = "'219"   # Max Forecast
$wsSummary.Cells.Item(14, 2).Value = "'69"    # Min Forecast
